$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1744186046511628
$ws.Range("C2").Value = 0.5736434108527132
$ws.Range("J2").Value = 0.01937984496124031
$ws.Range("P2").Value = 0.1085271317829457
$ws.Range("S2").Value = 0.124031007751938
$ws.Range("B3").Value = 0.006493506493506494
$ws.Range("C3").Value = 0.05194805194805195
$ws.Range("J3").Value = 0.02597402597402598
$ws.Range("P3").Value = 0.6623376623376623
$ws.Range("S3").Value = 0.2532467532467532
$ws.Range("P4").Value = 0.65
$ws.Range("S4").Value = 0.35
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.06091370558375635
$ws.Range("D6").Value = 0.01015228426395939
$ws.Range("F6").Value = 0.09644670050761421
$ws.Range("J6").Value = 0.2284263959390863
$ws.Range("O6").Value = 0.02030456852791878
$ws.Range("Q6").Value = 0.1116751269035533
$ws.Range("R6").Value = 0.06091370558375635
$ws.Range("S6").Value = 0.4111675126903553
$ws.Range("B7").Value = 0.08743169398907104
$ws.Range("D7").Value = 0.0273224043715847
$ws.Range("F7").Value = 0.06010928961748634
$ws.Range("J7").Value = 0.1256830601092896
$ws.Range("O7").Value = 0.01639344262295082
$ws.Range("Q7").Value = 0.1584699453551913
$ws.Range("R7").Value = 0.09289617486338798
$ws.Range("S7").Value = 0.4316939890710382
$ws.Range("B8").Value = 0.1075
$ws.Range("D8").Value = 0.01
$ws.Range("E8").Value = 0.0025
$ws.Range("F8").Value = 0.0575
$ws.Range("J8").Value = 0.105
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.1525
$ws.Range("R8").Value = 0.095
$ws.Range("S8").Value = 0.45
$ws.Range("B9").Value = 0.1428571428571428
$ws.Range("D9").Value = 0.02040816326530612
$ws.Range("E9").Value = 0.006802721088435374
$ws.Range("F9").Value = 0.06122448979591837
$ws.Range("J9").Value = 0.04761904761904762
$ws.Range("O9").Value = 0.01360544217687075
$ws.Range("Q9").Value = 0.1972789115646258
$ws.Range("R9").Value = 0.08163265306122448
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.1149312377210216
$ws.Range("D10").Value = 0.0275049115913556
$ws.Range("E10").Value = 0.002946954813359528
$ws.Range("F10").Value = 0.06777996070726916
$ws.Range("J10").Value = 0.1159135559921415
$ws.Range("O10").Value = 0.02161100196463654
$ws.Range("Q10").Value = 0.1836935166994106
$ws.Range("R10").Value = 0.07269155206286837
$ws.Range("S10").Value = 0.3929273084479371
$ws.Range("G11").Value = 0.1661129568106312
$ws.Range("J11").Value = 0.1096345514950166
$ws.Range("K11").Value = 0.2093023255813954
$ws.Range("L11").Value = 0.5016611295681063
$ws.Range("S11").Value = 0.0132890365448505
$ws.Range("G12").Value = 0.7350993377483444
$ws.Range("J12").Value = 0.2185430463576159
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.01986754966887417
$ws.Range("S12").Value = 0.01986754966887417
$ws.Range("G13").Value = 0.7435897435897436
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.02564102564102564
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.009345794392523364
$ws.Range("H15").Value = 0.1728971962616822
$ws.Range("I15").Value = 0.06542056074766354
$ws.Range("J15").Value = 0.3364485981308411
$ws.Range("K15").Value = 0.0514018691588785
$ws.Range("M15").Value = 0.02803738317757009
$ws.Range("O15").Value = 0.102803738317757
$ws.Range("S15").Value = 0.2336448598130841
$ws.Range("F16").Value = 0.01910828025477707
$ws.Range("H16").Value = 0.2101910828025478
$ws.Range("I16").Value = 0.07643312101910828
$ws.Range("J16").Value = 0.3885350318471338
$ws.Range("K16").Value = 0.1082802547770701
$ws.Range("O16").Value = 0.04458598726114649
$ws.Range("S16").Value = 0.1528662420382166
$ws.Range("F17").Value = 0.02439024390243903
$ws.Range("H17").Value = 0.2073170731707317
$ws.Range("I17").Value = 0.07317073170731707
$ws.Range("J17").Value = 0.3201219512195122
$ws.Range("K17").Value = 0.125
$ws.Range("M17").Value = 0.01829268292682927
$ws.Range("N17").Value = 0.003048780487804878
$ws.Range("O17").Value = 0.08536585365853659
$ws.Range("S17").Value = 0.1432926829268293
$ws.Range("F18").Value = 0.03267973856209151
$ws.Range("H18").Value = 0.1699346405228758
$ws.Range("I18").Value = 0.09803921568627451
$ws.Range("J18").Value = 0.4313725490196079
$ws.Range("K18").Value = 0.130718954248366
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.07843137254901961
$ws.Range("F19").Value = 0.02116402116402116
$ws.Range("H19").Value = 0.2107583774250441
$ws.Range("I19").Value = 0.07319223985890652
$ws.Range("J19").Value = 0.3606701940035273
$ws.Range("K19").Value = 0.1305114638447972
$ws.Range("M19").Value = 0.02557319223985891
$ws.Range("N19").Value = 0.001763668430335097
$ws.Range("O19").Value = 0.07671957671957672
$ws.Range("S19").Value = 0.09964726631393298
